$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark from the "Idea: ..." paragraph
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Add a new list paragraph after "Must be very careful with arrays..."
#    with the same ListParagraph / numId=3 formatting, and place the
#    _GoBack bookmark at the end of its text (collapsed range).
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
$newLast = $d.Paragraphs.Last

# Use a temporary trailing placeholder character so the bookmark insertion
# point is not the absolute last position in the document (a boundary case
# that the COM host mishandles for zero-length bookmark ranges).
$newLast.Range.Text = "Ohmergerd, remember to change the baudrate from 9600 to 31250 when going from printing to Serial to MIDI communicationX"

$pos = $d.Content.End - 2
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the temporary placeholder character
$xRange = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$xRange.Delete()
